$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 126
$ws.Range("A126").Value = "(홍보)KSC2025 카카오X한국정보과학회 AI 에이전트 경진대회(12.17, 여수엑스포컨벤션센터, 서류접수 ~11/18)"
$ws.Range("B126").Value = "공주대학교SW중심대학사업단"
$ws.Range("C126").NumberFormat = "@"
$ws.Range("C126").Value = "2025-11-17"
$ws.Range("C126").NumberFormat = "yyyy\-mm"
$ws.Range("D126").Value = "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=139"

# Row 127
$ws.Range("A127").Value = "2025년 SW알고리즘 경진대회 수상자 안내"
$ws.Range("B127").Value = "공주대학교SW중심대학사업단"
$ws.Range("C127").NumberFormat = "@"
$ws.Range("C127").Value = "2025-11-17"
$ws.Range("C127").NumberFormat = "yyyy\-mm"
$ws.Range("D127").Value = "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=138"
